$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the month header (row 5, column A)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2 = "شهریور، مهر و آبان 99"

# ---------------------------------------------------------------------------
# 2. Update the existing activity rows (6,7,8) - text / hours / notes
# ---------------------------------------------------------------------------
$ws.Range("B6").Value2 = "* Code Refactoring"
$ws.Range("C6").Value2 = 6
$ws.Range("E6").Value2 = "• Adjusted orientation widget for the phantom"

$ws.Range("C7").Value2 = 8
$ws.Range("E7").Value2 = "• Adjusted 2D views for the phantom"

$ws.Range("C8").Value2 = 14
# E8 previously empty - copy the format used by the other "notes" cells (E7)
$ws.Range("E7").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value2 = "• Extract/Load Image Centerline"

# ---------------------------------------------------------------------------
# 3. Insert two brand-new activity rows after row 8 (new rows 9 & 10)
# ---------------------------------------------------------------------------
$ws.Rows("9:10").Insert()

# Copy formatting from row 8 (B/C) and row 7 (E) into the freshly inserted rows
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Range("B8").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("B9").Value2 = "* Tracker"
$ws.Range("C9").Value2 = 2
$ws.Range("E9").Value2 = "• Record/Load Tracker Centerline"

$ws.Range("B10").Value2 = "* 2D/3D Views"
$ws.Range("C10").Value2 = 2
$ws.Range("E10").Value2 = "• Integrated Registration Process"

# ---------------------------------------------------------------------------
# 4. Row that used to be "* Online Tracking" (old row 9) is now row 11 and
#    becomes "* Patients Database"
# ---------------------------------------------------------------------------
$ws.Range("B11").Value2 = "* Patients Database"
$ws.Range("C11").Value2 = 1

# ---------------------------------------------------------------------------
# 5. Row that used to be "* Meetings & other" (old row 10) is now row 12.
#    It also inherits the taller font / row height that used to sit on the
#    old Total row (old row 11, E11 s="18").
# ---------------------------------------------------------------------------
$ws.Range("C12").Value2 = 2

$ws.Range("E13").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E13").Clear()
$ws.Rows(12).RowHeight = 15.6

# ---------------------------------------------------------------------------
# 6. Totals / calculated rows shift from 11-13 to 13-15; formulas are
#    automatically re-targeted by the row insert, just refresh the literal
#    "non billable" hours value.
# ---------------------------------------------------------------------------
$ws.Range("D14").Value2 = 0

# ---------------------------------------------------------------------------
# 7. Column A width changed (18.33 -> 21)
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 20.17

# ---------------------------------------------------------------------------
# 8. Selection cosmetics (matches the saved file's last selected cell)
# ---------------------------------------------------------------------------
$ws.Range("E18").Select()
